$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Fill in a new work-hours entry on row 17: date, hours, and location note.
# Match the date formatting already used by the rows above (e.g. A16).
$ws.Cells.Item(16, 1).Copy() | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(17, 1).Value = 41408
$ws.Cells.Item(17, 2).Value = 5
$ws.Cells.Item(17, 3).Value = "projektihuoneella"

# Move the active selection to the recommended next empty row (C19).
$ws.Range("C19").Select()
